$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 342.8889
$ws.Range("I8").Value = 163.625
$ws.Range("J8").Value = 1777
$ws.Range("K8").Value = 490.875
$ws.Range("L8").Value = 5331
$ws.Range("M8").Value = -351.875
$ws.Range("N8").Value = -5609
$ws.Range("H9").Value = 128.4
$ws.Range("I9").Value = 110.5
$ws.Range("K9").Value = 110.5
$ws.Range("M9").Value = 58.5
$ws.Range("H19").Value = 438.82144
$ws.Range("I19").Value = 341.6154
$ws.Range("J19").Value = 523.06665
$ws.Range("K19").Value = 341.6154
$ws.Range("L19").Value = 523.06665
$ws.Range("M19").Value = -166.6154
$ws.Range("N19").Value = -873.06665
$ws.Range("H28").Value = 292.52942
$ws.Range("I28").Value = 212.25
$ws.Range("K28").Value = 212.25
$ws.Range("M28").Value = 272.75
$ws.Range("H88").Value = 8978.579
$ws.Range("I88").Value = 2899
$ws.Range("J88").Value = 10118.5
$ws.Range("K88").Value = 2899
$ws.Range("L88").Value = 10118.5
$ws.Range("M88").Value = -2493
$ws.Range("N88").Value = -10930.5
$ws.Range("H91").Value = 8978.579
$ws.Range("I91").Value = 2899
$ws.Range("J91").Value = 10118.5
$ws.Range("K91").Value = 2899
$ws.Range("L91").Value = 10118.5
$ws.Range("M91").Value = -1495
$ws.Range("N91").Value = -12926.5
$ws.Range("H132").Value = 6668641.5
$ws.Range("I132").Value = 1296.5952
$ws.Range("J132").Value = 41672204
$ws.Range("K132").Value = 3889.7856
$ws.Range("L132").Value = 125016612
$ws.Range("M132").Value = -1359.7856
$ws.Range("N132").Value = -125021672
$ws.Range("H141").Value = 2134.077
$ws.Range("I141").Value = 1222.0714
$ws.Range("K141").Value = 3666.2142
$ws.Range("M141").Value = 1513.7858
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 167661.28
$ws.Range("I61").Value = 3927.3022
$ws.Range("J61").Value = 558803.5600000001
$ws.Range("K61").Value = 3927.3022
$ws.Range("L61").Value = 558803.5600000001
$ws.Range("M61").Value = -3715.3022
$ws.Range("N61").Value = -559227.5600000001
$ws.Range("H97").Value = 1586.9412
$ws.Range("I97").Value = 1462.7142
$ws.Range("J97").Value = 2166.6667
$ws.Range("K97").Value = 1462.7142
$ws.Range("L97").Value = 2166.6667
$ws.Range("M97").Value = -966.7141999999999
$ws.Range("N97").Value = -3158.6667
$ws.Range("H122").Value = 586213.75
$ws.Range("I122").Value = 736263
$ws.Range("J122").Value = 2688.889
$ws.Range("K122").Value = 2208789
$ws.Range("L122").Value = 8066.667
$ws.Range("M122").Value = -2206339
$ws.Range("N122").Value = -12966.667
$ws.Range("H136").Value = 167661.28
$ws.Range("I136").Value = 3927.3022
$ws.Range("J136").Value = 558803.5600000001
$ws.Range("K136").Value = 11781.9066
$ws.Range("L136").Value = 1676410.68
$ws.Range("M136").Value = -9231.9066
$ws.Range("N136").Value = -1681510.68
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8253.263000000001
$ws.Range("I20").Value = 1676.6
$ws.Range("J20").Value = 20900.691
$ws.Range("K20").Value = 1676.6
$ws.Range("L20").Value = 20900.691
$ws.Range("M20").Value = -1429.6
$ws.Range("N20").Value = -21394.691
$ws.Range("H107").Value = 233219.92
$ws.Range("I107").Value = 302675.9
$ws.Range("J107").Value = 1700
$ws.Range("K107").Value = 302675.9
$ws.Range("L107").Value = 1700
$ws.Range("M107").Value = -300755.9
$ws.Range("N107").Value = -5540
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10970764
$ws.Range("I31").Value = 1703.5186
$ws.Range("J31").Value = 25778996
$ws.Range("K31").Value = 1703.5186
$ws.Range("L31").Value = 25778996
$ws.Range("M31").Value = -1408.5186
$ws.Range("N31").Value = -25779586
$ws.Range("H34").Value = 10970764
$ws.Range("I34").Value = 1703.5186
$ws.Range("J34").Value = 25778996
$ws.Range("K34").Value = 1703.5186
$ws.Range("L34").Value = 25778996
$ws.Range("M34").Value = -1501.5186
$ws.Range("N34").Value = -25779400
$ws.Range("H134").Value = 7307499.5
$ws.Range("I134").Value = 11497299
$ws.Range("J134").Value = 557266.75
$ws.Range("K134").Value = 34491897
$ws.Range("L134").Value = 1671800.25
$ws.Range("M134").Value = -34489362
$ws.Range("N134").Value = -1676870.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 6250
$ws.Range("J55").Value = 6250
$ws.Range("L55").Value = 18750
$ws.Range("N55").Value = -19104
$ws.Range("H107").Value = 454.5
$ws.Range("I107").Value = 420.33334
$ws.Range("J107").Value = 501.0909
$ws.Range("K107").Value = 1261.00002
$ws.Range("L107").Value = 1503.2727
$ws.Range("M107").Value = 658.9999800000001
$ws.Range("N107").Value = -5343.2727
$ws.Range("H129").Value = 1335.8077
$ws.Range("I129").Value = 1280
$ws.Range("J129").Value = 1370.6875
$ws.Range("K129").Value = 3840
$ws.Range("L129").Value = 4112.0625
$ws.Range("M129").Value = 1160
$ws.Range("N129").Value = -14112.0625
$ws.Range("H133").Value = 28352.17
$ws.Range("I133").Value = 104871.9
$ws.Range("J133").Value = 7671.162
$ws.Range("K133").Value = 314615.7
$ws.Range("L133").Value = 23013.486
$ws.Range("M133").Value = -309555.7
$ws.Range("N133").Value = -33133.486
$ws.Range("H136").Value = 6270.4414
$ws.Range("I136").Value = 9266.25
$ws.Range("J136").Value = 4636.364
$ws.Range("K136").Value = 27798.75
$ws.Range("L136").Value = 13909.092
$ws.Range("M136").Value = -22698.75
$ws.Range("N136").Value = -24109.092
$ws.Range("H137").Value = 16618.916
$ws.Range("I137").Value = 10370.77
$ws.Range("J137").Value = 24003.092
$ws.Range("K137").Value = 31112.31
$ws.Range("L137").Value = 72009.276
$ws.Range("M137").Value = -26012.31
$ws.Range("N137").Value = -82209.276
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2797.6667
$ws.Range("I102").Value = 2668.4285
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 2668.4285
$ws.Range("L102").Value = 3250
$ws.Range("M102").Value = -1046.4285
$ws.Range("N102").Value = -6494
$ws.Range("H122").Value = 113737656
$ws.Range("I122").Value = 212964560
$ws.Range("J122").Value = 42861290
$ws.Range("K122").Value = 638893680
$ws.Range("L122").Value = 128583870
$ws.Range("M122").Value = -638891230
$ws.Range("N122").Value = -128588770
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 50002090
$ws.Range("I68").Value = 1954
$ws.Range("J68").Value = 100002220
$ws.Range("K68").Value = 1954
$ws.Range("L68").Value = 100002220
$ws.Range("M68").Value = -1205
$ws.Range("N68").Value = -100003718
$ws.Range("H71").Value = 50002090
$ws.Range("I71").Value = 1954
$ws.Range("J71").Value = 100002220
$ws.Range("K71").Value = 9770
$ws.Range("L71").Value = 500011100
$ws.Range("M71").Value = -6026
$ws.Range("N71").Value = -500018588
$ws.Range("H136").Value = 9754
$ws.Range("I136").Value = 7010.409
$ws.Range("J136").Value = 15789.9
$ws.Range("K136").Value = 21031.227
$ws.Range("L136").Value = 47369.7
$ws.Range("M136").Value = -18481.227
$ws.Range("N136").Value = -52469.7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 76923816
$ws.Range("I107").Value = 200000560
$ws.Range("J107").Value = 856.25
$ws.Range("K107").Value = 600001680
$ws.Range("L107").Value = 2568.75
$ws.Range("M107").Value = -599999760
$ws.Range("N107").Value = -6408.75
$ws.Range("H122").Value = 833.8333
$ws.Range("I122").Value = 833.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2501.4999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -51.4998999999998
$ws.Range("N122").ClearContents()
